# Adds a new "2022-Q3" worksheet (with fund holdings data) right after
# the "总计" (summary) sheet, and updates the "总计" sheet to add a new
# leading row for 2022-Q3, shifting all other quarters down by one row.

$wb = $excel.ActiveWorkbook

# --- Locate the "总计" (summary) sheet; it stays sheet #1 ---
$summary = $wb.Worksheets.Item(1)

# --- Create the new "2022-Q3" worksheet right after "总计" ---
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# --- Populate the header row ---
$q3.Cells.Item(1,2).Value = "基金代码"
$q3.Cells.Item(1,3).Value = "基金名称"
$q3.Cells.Item(1,4).Value = "基金规模"
$q3.Cells.Item(1,5).Value = "股票总仓位"
$q3.Cells.Item(1,6).Value = "仓位占比"
$q3.Cells.Item(1,7).Value = "持有市值(亿元)"
$q3.Cells.Item(1,8).Value = "仓位排名"

# --- Populate the data rows (2..15) ---
# Numeric-looking text fields (fund code, scale, position %, etc.) are
# entered with a leading apostrophe so Excel stores them as text (not
# auto-converted numbers), matching the source data which keeps values
# like fund codes ("004350") and formatted numbers ("0.1370") as text.
$q3.Cells.Item(2,1).Value = 0
$q3.Cells.Item(2,2).Value = "'004350"
$q3.Cells.Item(2,3).Value = "汇丰晋信价值先锋股票A"
$q3.Cells.Item(2,4).Value = "'5.31"
$q3.Cells.Item(2,5).Value = "'94.44"
$q3.Cells.Item(2,6).Value = "'2.58"
$q3.Cells.Item(2,7).Value = "'0.1370"
$q3.Cells.Item(2,8).Value = 8
$q3.Cells.Item(3,1).Value = 1
$q3.Cells.Item(3,2).Value = "'700001"
$q3.Cells.Item(3,3).Value = "平安行业先锋混合"
$q3.Cells.Item(3,4).Value = "'1.82"
$q3.Cells.Item(3,5).Value = "'91.01"
$q3.Cells.Item(3,6).Value = "'3.48"
$q3.Cells.Item(3,7).Value = "'0.0633"
$q3.Cells.Item(3,8).Value = 6
$q3.Cells.Item(4,1).Value = 2
$q3.Cells.Item(4,2).Value = "'290012"
$q3.Cells.Item(4,3).Value = "泰信行业精选灵活配置混合A"
$q3.Cells.Item(4,4).Value = "'0.75"
$q3.Cells.Item(4,5).Value = "'91.96"
$q3.Cells.Item(4,6).Value = "'5.93"
$q3.Cells.Item(4,7).Value = "'0.0445"
$q3.Cells.Item(4,8).Value = 3
$q3.Cells.Item(5,1).Value = 3
$q3.Cells.Item(5,2).Value = "'013767"
$q3.Cells.Item(5,3).Value = "平安价值回报混合A"
$q3.Cells.Item(5,4).Value = "'0.88"
$q3.Cells.Item(5,5).Value = "'92.26"
$q3.Cells.Item(5,6).Value = "'3.28"
$q3.Cells.Item(5,7).Value = "'0.0289"
$q3.Cells.Item(5,8).Value = 5
$q3.Cells.Item(6,1).Value = 4
$q3.Cells.Item(6,2).Value = "'003132"
$q3.Cells.Item(6,3).Value = "德邦新回报灵活配置混合"
$q3.Cells.Item(6,4).Value = "'0.62"
$q3.Cells.Item(6,5).Value = "'72.03"
$q3.Cells.Item(6,6).Value = "'2.51"
$q3.Cells.Item(6,7).Value = "'0.0156"
$q3.Cells.Item(6,8).Value = 8
$q3.Cells.Item(7,1).Value = 5
$q3.Cells.Item(7,2).Value = "'001900"
$q3.Cells.Item(7,3).Value = "诺安精选价值混合"
$q3.Cells.Item(7,4).Value = "'0.12"
$q3.Cells.Item(7,5).Value = "'85.83"
$q3.Cells.Item(7,6).Value = "'3.00"
$q3.Cells.Item(7,7).Value = "'0.0036"
$q3.Cells.Item(7,8).Value = 6
$q3.Cells.Item(8,1).Value = 6
$q3.Cells.Item(8,2).Value = "'001664"
$q3.Cells.Item(8,3).Value = "平安鑫安混合A"
$q3.Cells.Item(8,4).Value = "'0.33"
$q3.Cells.Item(8,5).Value = "'27.84"
$q3.Cells.Item(8,6).Value = "'1.03"
$q3.Cells.Item(8,7).Value = "'0.0034"
$q3.Cells.Item(8,8).Value = 6
$q3.Cells.Item(9,1).Value = 7
$q3.Cells.Item(9,2).Value = "'006433"
$q3.Cells.Item(9,3).Value = "平安鑫利灵活配置混合C"
$q3.Cells.Item(9,4).Value = "'0.23"
$q3.Cells.Item(9,5).Value = "'27.39"
$q3.Cells.Item(9,6).Value = "'1.03"
$q3.Cells.Item(9,7).Value = "'0.0024"
$q3.Cells.Item(9,8).Value = 6
$q3.Cells.Item(10,1).Value = 8
$q3.Cells.Item(10,2).Value = "'002583"
$q3.Cells.Item(10,3).Value = "泰信行业精选灵活配置混合C"
$q3.Cells.Item(10,4).Value = "'0.04"
$q3.Cells.Item(10,5).Value = "'91.96"
$q3.Cells.Item(10,6).Value = "'5.93"
$q3.Cells.Item(10,7).Value = "'0.0024"
$q3.Cells.Item(10,8).Value = 3
$q3.Cells.Item(11,1).Value = 9
$q3.Cells.Item(11,2).Value = "'007049"
$q3.Cells.Item(11,3).Value = "平安鑫安混合E"
$q3.Cells.Item(11,4).Value = "'0.16"
$q3.Cells.Item(11,5).Value = "'27.84"
$q3.Cells.Item(11,6).Value = "'1.03"
$q3.Cells.Item(11,7).Value = "'0.0016"
$q3.Cells.Item(11,8).Value = 6
$q3.Cells.Item(12,1).Value = 10
$q3.Cells.Item(12,2).Value = "'003626"
$q3.Cells.Item(12,3).Value = "平安鑫利灵活配置混合A"
$q3.Cells.Item(12,4).Value = "'0.15"
$q3.Cells.Item(12,5).Value = "'27.39"
$q3.Cells.Item(12,6).Value = "'1.03"
$q3.Cells.Item(12,7).Value = "'0.0015"
$q3.Cells.Item(12,8).Value = 6
$q3.Cells.Item(13,1).Value = 11
$q3.Cells.Item(13,2).Value = "'013768"
$q3.Cells.Item(13,3).Value = "平安价值回报混合C"
$q3.Cells.Item(13,4).Value = "'0.04"
$q3.Cells.Item(13,5).Value = "'92.26"
$q3.Cells.Item(13,6).Value = "'3.28"
$q3.Cells.Item(13,7).Value = "'0.0013"
$q3.Cells.Item(13,8).Value = 5
$q3.Cells.Item(14,1).Value = 12
$q3.Cells.Item(14,2).Value = "'015364"
$q3.Cells.Item(14,3).Value = "汇丰晋信价值先锋股票C"
$q3.Cells.Item(14,4).Value = "'0.02"
$q3.Cells.Item(14,5).Value = "'94.44"
$q3.Cells.Item(14,6).Value = "'2.58"
$q3.Cells.Item(14,7).Value = "'0.0005"
$q3.Cells.Item(14,8).Value = 8
$q3.Cells.Item(15,1).Value = 13
$q3.Cells.Item(15,2).Value = "'001665"
$q3.Cells.Item(15,3).Value = "平安鑫安混合C"
$q3.Cells.Item(15,4).Value = "'0.01"
$q3.Cells.Item(15,5).Value = "'27.84"
$q3.Cells.Item(15,6).Value = "'1.03"
$q3.Cells.Item(15,7).Value = "'0.0001"
$q3.Cells.Item(15,8).Value = 6

# --- Apply the bold/bordered "header" style (same as other quarter sheets) ---
# Copy formats from the summary sheet's styled cells: B1:D1 (header) and A2 (index column)
$summary.Range("B1:D1").Copy() | Out-Null
$q3.Range("B1:H1").PasteSpecial(-4122) | Out-Null
$summary.Range("A2").Copy() | Out-Null
$q3.Range("A2:A15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Update the "总计" summary sheet: insert a new leading row for 2022-Q3 ---
$summary.Rows.Item(2).Insert()
# Row-insert in Excel copies the formatting of the row above (the header),
# which is not what we want for a plain data row, so clear it back to default.
$summary.Range("A2:D2").ClearFormats()

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 14
$summary.Cells.Item(2,4).Value = 0.31
$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(3,2).Value = "2022-Q2"
$summary.Cells.Item(3,3).Value = 27
$summary.Cells.Item(3,4).Value = 3.68
$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(4,2).Value = "2022-Q1"
$summary.Cells.Item(4,3).Value = 26
$summary.Cells.Item(4,4).Value = 5.07
$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(5,2).Value = "2021-Q4"
$summary.Cells.Item(5,3).Value = 5
$summary.Cells.Item(5,4).Value = 2.64
$summary.Cells.Item(6,1).Value = 4
$summary.Cells.Item(6,2).Value = "2021-Q3"
$summary.Cells.Item(6,3).Value = 3
$summary.Cells.Item(6,4).Value = 2.59
$summary.Cells.Item(7,1).Value = 5
$summary.Cells.Item(7,2).Value = "2021-Q2"
$summary.Cells.Item(7,3).Value = 2
$summary.Cells.Item(7,4).Value = 2.68
$summary.Cells.Item(8,1).Value = 6
$summary.Cells.Item(8,2).Value = "2021-Q1"
$summary.Cells.Item(8,3).Value = 11
$summary.Cells.Item(8,4).Value = 4.1

# --- Re-apply the bold/bordered style to A2 of the summary sheet (the index column) ---
$summary.Range("A3").Copy() | Out-Null
$summary.Range("A2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

Write-Host "2022-Q3 sheet added and 总计 summary updated"
